$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $ok = $d.Content.Find.Execute(
        $findText, $true, $false, $false, $false, $false,
        $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $findText"
    }
}

# 1) Intro paragraph: drop "Nu's scripting system as well " before "MetaFunctions".
Replace-Text `
    "the semantic design for Nu’s scripting system as well MetaFunctions" `
    "the semantic design for MetaFunctions"

# 2) Second-capability paragraph: "intended semantics" -> "intended structure".
Replace-Text `
    "specify their program’s intended semantics in one of two ways" `
    "specify their program’s intended structure in one of two ways"

# 3) Denotational-design bullet: "whose semantics can be specified" -> "whose
#    structure can be specified".
Replace-Text `
    "programs / subprograms whose semantics can be specified" `
    "programs / subprograms whose structure can be specified"

# 4) Semantic-design bullet: drop the word "semantic" from both occurrences of
#    "level of semantic detail".
Replace-Text `
    "for the level of semantic detail at which" `
    "for the level of detail at which"

Replace-Text `
    "increase the level of semantic detail for designs" `
    "increase the level of detail for designs"

# 5) Register the ListLabel46 character style that the saved package also
#    gains (same shape as the existing ListLabel1..45 family: Courier New,
#    9pt, English).
$s = $d.Styles.Add("ListLabel 46", 2)
$s.Font.NameAscii = "Courier New"
$s.Font.Name = "Courier New"
$s.Font.Size = 9
$s.LanguageID = "en"
$s.QuickStyle = $true
